$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.621.44"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "3.089.83"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("D9").Value = "3.082.00"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "3.601.90"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "63.573.23"
$ws.Range("D19").Value = "3.085.36"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +4.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("D35").Value = "0.0₃0852"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  -2.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "444.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.284"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0363"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("D47").Value = "2.803.09"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.05%  "
